$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37, pushing the old rows 37-39 down to 38-40.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new weekly price record.
$ws.Cells.Item(37, 1).Value = 1
$ws.Cells.Item(37, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(37, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(37, 4).Value = 45008
$ws.Cells.Item(37, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(37, 5).Value = 15
$ws.Cells.Item(37, 6).Value = 100112003
$ws.Cells.Item(37, 7).Value = "Ajo"
$ws.Cells.Item(37, 8).Value = "Chino"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 1750
$ws.Cells.Item(37, 11).Value = 18000
$ws.Cells.Item(37, 12).Value = 19000
$ws.Cells.Item(37, 13).Value = 18500
$ws.Cells.Item(37, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(37, 15).Value = "China"
$ws.Cells.Item(37, 16).Value = 1850
$ws.Cells.Item(37, 17).Value = 10
$ws.Cells.Item(37, 18).Value = "Hortaliza"
